# Realizando envio do relatorio da busca das palavras chave por e-mail
#
# Adds a title/subtitle banner ("Usiwal" / "usiwal.com.br") above the
# existing "Palavra Chave / Pagina / Posicao" header, and fills in the
# keyword search-result rows underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 1 (Palavra Chave / Página / Posição) becomes row 4, so clear
# the two cells we are not about to overwrite with new text before we turn
# them into the styled-but-empty B1/C1 of the new title row.
$ws.Range("B1").ClearContents() | Out-Null
$ws.Range("C1").ClearContents() | Out-Null

# --- Row 1: big bold title "Usiwal" --------------------------------------
$ws.Range("A1").Value = "Usiwal"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Font.Name = "Arial Black"
$ws.Range("A1:C1").Font.Size = 24
$ws.Range("A1:C1").Font.Color = 3875864   # RGB(0x18, 0x24, 0x3B) == #18243B
$ws.Range("A1:C1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:C1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1:C1").Merge() | Out-Null

# --- Row 2: subtitle "usiwal.com.br" --------------------------------------
$ws.Range("A2").Value = "usiwal.com.br"
$ws.Range("A2:C2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A2:C2").VerticalAlignment = -4108     # xlCenter
$ws.Range("A2:C2").Merge() | Out-Null

# --- Row 4: table header (kept from the original sheet) -------------------
$ws.Range("A4").Value = "Palavra Chave"
$ws.Range("B4").Value = "Página"
$ws.Range("C4").Value = "Posição"

# --- Rows 5-6: keyword search results --------------------------------------
$ws.Range("A5").Value = "Cilindros Hidráulicos De Alta Pressão"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

$ws.Range("A6").Value = "Cilindros Hidráulicos 700 Bar"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1

# --- Column widths ----------------------------------------------------------
# ColumnWidth is in "characters"; the engine adds ~5/6 of a character of
# padding when it stores the OOXML <col width>, so back that out here so the
# saved width lands on exactly 100 / 25 characters.
$ws.Columns.Item(1).ColumnWidth = 99.16666666666667
$ws.Columns.Item(2).ColumnWidth = 24.166666666666668
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668
